$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Status text for "CRUD Branch" (row 2, column F)
$ws.Range("F2").Value = "CU - Complete`nR - Pending"

# Add new Status text for "CRUD Agent" (row 3, column F)
$ws.Range("F3").Value = "CU - SP Done; Calling Pending`nR - Pending"

# Adjust row heights to match new wrapped text content
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 43.2

# Update selection to F3
$ws.Range("F3").Select()
